$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cell A14 (previously held the stray "Q" label)
$ws.Range("A14").Value = ""

# Rename attribute labels for better clarification
$ws.Range("B27").Value = "country_name"
$ws.Range("B29").Value = "org_name"
$ws.Range("B30").Value = "org_type"
$ws.Range("B28").Value = "etc_code"
$ws.Range("D28").Value = "Ebola Treatment Center / Ebola Care Facility's identification code"

# Update the selected cell to match the saved cursor position
$ws.Range("D28").Select()
